$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D9").NumberFormat = "@"
$ws.Range("D2:D9").Value = "'TRUE"

$ws.Range("D2:D9").Select()
